$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 96; existing rows 96..124 shift down to 97..125
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with its data
$ws.Range("A96").Value = 5
$ws.Range("B96").Value = "Macroferia Regional de Talca"
$ws.Range("C96").Value = "Maule"
$ws.Range("D96").Value = 44900
$ws.Range("E96").Value = 7
$ws.Range("F96").Value = 100112022
$ws.Range("G96").Value = "Arveja Verde"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 400
$ws.Range("K96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("M96").Value = 20000
$ws.Range("N96").Value = "`$/saco 25 kilos"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 800
$ws.Range("Q96").Value = 25
$ws.Range("R96").Value = "Hortaliza"
